$d = $word.ActiveDocument

# Several "newly downloaded" <id>...</id> tags are currently split across
# three runs: the literal "<id>" text (Courier New / gold), the bare id
# value (plain black), and the literal "</id>" text (Courier New / gold
# again). They should instead live in a single run/string, i.e.
# "<id>p145r_2</id>" etc., formatted like the surrounding tag markup.
#
# For each id, locate the full "<id>VALUE</id>" span, trim the match
# down to just the leading "<id>" (still inside the first run), delete
# everything after it up to the end of the match, and then re-type the
# removed "VALUE</id>" text right back in immediately afterwards. Typing
# directly after the end of the existing run's text appends to that run
# instead of minting new ones, so the three runs collapse into the one
# run that originally held "<id>".

$ids = @("p145r_2", "p145v_1", "p145v_2")

foreach ($id in $ids) {
    $needle = "<id>" + $id + "</id>"

    $match = $d.Content
    $found = $match.Find.Execute($needle, $false, $true, $false, $false, $false, `
                                  $true, 1, $false, $null, 0)

    if ($found) {
        $tagStart = $match.Start
        $tagEnd = $match.End
        $openTagEnd = $tagStart + 4   # length of "<id>"

        $rest = $d.Range($openTagEnd, $tagEnd)
        $rest.Delete()

        $insertionPoint = $d.Range($openTagEnd, $openTagEnd)
        $insertionPoint.InsertAfter($id + "</id>")
    }
}
